$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 23.218118119056388
$ws.Range("C2").Value = 46.915574224234724
$ws.Range("D2").Value = 20.821599490754718
$ws.Range("E2").Value = 44.053825463366024

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 21.561630270302388
$ws.Range("C3").Value = 41.479346961995333
$ws.Range("D3").Value = 32.706277011313034
$ws.Range("E3").Value = 28.425872578541874

# Update the selection to match the new committed range
$ws.Range("B1:E3").Select()
